# Insert a new weekly data row for "Cebolla" (Hortaliza, Macroferia Regional de Talca)
# above the existing row 401. This pushes rows 401:507 down to 402:508 and
# extends the sheet's used range to A1:R508.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("401:401").Insert()

$ws.Range("A401").Value = 5
$ws.Range("B401").Value = "Macroferia Regional de Talca"
$ws.Range("C401").Value = "Maule"
$ws.Range("D401").Value = 44642
$ws.Range("E401").Value = 7
$ws.Range("F401").Value = 100112004
$ws.Range("G401").Value = "Cebolla"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "1a (cosecha)"
$ws.Range("J401").Value = 2500
$ws.Range("K401").Value = 4500
$ws.Range("L401").Value = 4500
$ws.Range("M401").Value = 4500
$ws.Range("N401").Value = "`$/malla 25 kilos"
$ws.Range("O401").Value = "Región del Maule"
$ws.Range("P401").Value = 180
$ws.Range("Q401").Value = 25
$ws.Range("R401").Value = "Hortaliza"
